# Tweak personal statement grammar
#
# Rewrites the personal-statement paragraph ("I'm an experienced
# engineer, ...") to fix capitalisation, tighten the wording, and
# correct the "paassionate" typo, per the commit's intent.

$d = $word.ActiveDocument

# 1) "experienced engineer" -> "experienced Engineer"
$d.Content.Find.Execute(
    "experienced engineer,", $true, $false, $false, $false, $false,
    $true, 1, $false, "experienced Engineer,", 2) | Out-Null

# 2) Tighten "diverse projects. Specializes in optimizing platform
#    performance, ensuring security, and integrating services." into
#    "diverse projects, specialising in platform optimisation,
#    ensuring security, and integrating self-service."
$d.Content.Find.Execute(
    "diverse projects. Specializes in optimizing platform performance, ensuring security, and integrating services. I can",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "diverse projects, specialising in platform optimisation, ensuring security, and integrating self-service. I can",
    2) | Out-Null

# 3) Add a comma before "and I'm passionate" and fix the "paassionate" typo
$d.Content.Find.Execute(
    "observable solutions and I" + [char]0x2019 + "m paassionate",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "observable solutions, and I" + [char]0x2019 + "m passionate",
    2) | Out-Null
